$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("roundCount"), shifting
# isSystemValid..CT(B) one column to the right (C->D, D->E, ... T->U).
$ws.Range("C:C").Insert()

# Set the new header
$ws.Range("C1").Value = "roundCount"

# Row data for columns A..U (row 2 through row 7)
$data = @(
  @("SS512",128,20,20,20,20,0.006081915,0.008213257999999999,0.006675637,0.006746268,0.002869511,0.002893269,17366016,17366016,17407795.2,17408000,17408000,17408000,80,80,88),
  @("SS512",160,20,20,20,20,0.00625782,0.008428346999999999,0.006888354,0.006956375,0.002940571,0.003030062,17408000,17408000,17408000,17408000,17408000,17408000,80,80,88),
  @("SS512",224,20,20,20,20,0.006056547,0.008280503999999999,0.006794643,0.006845415,0.002927208,0.002963364,17408000,17408000,17408000,17408000,17408000,17408000,80,80,88),
  @("SS512",256,20,20,20,20,0.006344008,0.00874759,0.006993294,0.007171118,0.003034806,0.003063548,17408000,17408000,17408000,17408000,17408000,17408000,80,80,88),
  @("SS512",384,20,20,20,20,0.00606606,0.008221555,0.006723762,0.006720304,0.002907443,0.002882338,17408000,17408000,17408000,17408000,17408000,17408000,80,80,88),
  @("SS512",512,20,20,20,20,0.006339788,0.008527935,0.006845677,0.006907356,0.002990603,0.003015733,17408000,17408000,17408000,17408000,17408000,17408000,80,80,88)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}
